# Updated cryptos list with refreshed prices / 1h volume figures
# (also two pairs of rows were re-ordered: Polkadot/WrappedBTC and Stellar/Cosmos)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (NumberFormat "@") for every write so that
# numeric-looking strings (prices, percentages) are kept as literal
# text instead of being coerced into floating point numbers -
# matching the inlineStr cells in the original workbook. The style
# is reset back to Normal right after so no stray number format
# is left applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "28.348.99"
Set-TextValue $ws.Range("E2") "  -0.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.563.46"
Set-TextValue $ws.Range("E3") "  -0.04%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.09%  "

# Row 5
Set-TextValue $ws.Range("D5") "210.80"
Set-TextValue $ws.Range("E5") "  -0.40%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -0.61%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.11%  "

# Row 8
Set-TextValue $ws.Range("D8") "44.39"
Set-TextValue $ws.Range("E8") "  -4.30%  "

# Row 9
Set-TextValue $ws.Range("D9") "23.62"
Set-TextValue $ws.Range("E9") "  -1.98%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.244"
Set-TextValue $ws.Range("E10") "  -1.31%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0587"
Set-TextValue $ws.Range("E11") "  -0.94%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0893"
Set-TextValue $ws.Range("E12") "  +0.74%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.784.58"
Set-TextValue $ws.Range("E13") "  -0.14%  "

# Row 14
Set-TextValue $ws.Range("D14") "1.574.01"
Set-TextValue $ws.Range("E14") "  +0.88%  "

# Row 15
Set-TextValue $ws.Range("B15") "Polkadot"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "3.66"
Set-TextValue $ws.Range("E15") "  -0.59%  "

# Row 16
Set-TextValue $ws.Range("B16") "WrappedBTC"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "28.348.16"
Set-TextValue $ws.Range("E16") "  -0.56%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.513"
Set-TextValue $ws.Range("E17") "  -1.42%  "

# Row 18
Set-TextValue $ws.Range("D18") "61.04"
Set-TextValue $ws.Range("E18") "  -1.83%  "

# Row 19
Set-TextValue $ws.Range("D19") "228.06"
Set-TextValue $ws.Range("E19") "  -0.40%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.37"
Set-TextValue $ws.Range("E20") "  +0.46%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0679"
Set-TextValue $ws.Range("E21") "  -2.11%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -0.06%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +1.52%  "

# Row 24
Set-TextValue $ws.Range("D24") "8.93"
Set-TextValue $ws.Range("E24") "  -2.29%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -2.03%  "

# Row 26
Set-TextValue $ws.Range("D26") "150.37"
Set-TextValue $ws.Range("E26") "  -0.03%  "

# Row 27
Set-TextValue $ws.Range("D27") "14.90"
Set-TextValue $ws.Range("E27") "  -0.51%  "

# Row 28
Set-TextValue $ws.Range("B28") "Stellar"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D28") "0.103"
Set-TextValue $ws.Range("E28") "  +0.02%  "

# Row 29
Set-TextValue $ws.Range("B29") "Cosmos"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D29") "6.34"
Set-TextValue $ws.Range("E29") "  -1.54%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -0.13%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +2.16%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -4.05%  "

# Row 33
Set-TextValue $ws.Range("E33") "  -0.71%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.08"
Set-TextValue $ws.Range("E34") "  -0.72%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.385.57"
Set-TextValue $ws.Range("E35") "  -0.78%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +1.98%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.50"
Set-TextValue $ws.Range("E37") "  -3.14%  "

# Row 38
Set-TextValue $ws.Range("E38") "  -0.43%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +2.04%  "

# Row 40
Set-TextValue $ws.Range("E40") "  -1.95%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.520"
Set-TextValue $ws.Range("E41") "  -3.01%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +2.48%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -0.09%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.784"
Set-TextValue $ws.Range("E44") "  -0.58%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0470"
Set-TextValue $ws.Range("E45") "  -2.17%  "

# Row 46
Set-TextValue $ws.Range("D46") "5.34"
Set-TextValue $ws.Range("E46") "  -3.18%  "

# Row 47
Set-TextValue $ws.Range("D47") "62.30"
Set-TextValue $ws.Range("E47") "  -0.66%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.917"
Set-TextValue $ws.Range("E48") "  -6.15%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.697.48"
Set-TextValue $ws.Range("E49") "  -0.22%  "

# Row 50
Set-TextValue $ws.Range("D50") "85.46"
Set-TextValue $ws.Range("E50") "  -0.74%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0₆0100"
Set-TextValue $ws.Range("E51") "  -2.26%  "

Write-Output "Applied crypto list update."
